$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $escaped = $text.Replace('"', '""')
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

Set-TextValue 'D2' '65.561.48'
$ws.Range('E2').Value = '  -0.39%  '
Set-TextValue 'D3' '3.277.77'
$ws.Range('E3').Value = '  -0.68%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue 'D5' '578.15'
$ws.Range('E5').Value = '  +3.84%  '
Set-TextValue 'D6' '183.44'
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('E7').Value = '  +0.01%  '
Set-TextValue 'D8' '3.270.47'
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('E9').Value = '  -2.70%  '
Set-TextValue 'D10' '0.174'
$ws.Range('E10').Value = '  -5.46%  '
Set-TextValue 'D11' '0.569'
$ws.Range('E11').Value = '  -2.12%  '
Set-TextValue 'D12' '46.04'
$ws.Range('E12').Value = '  -3.17%  '
Set-TextValue 'D13' '0.0000261'
$ws.Range('E13').Value = '  -2.61%  '
Set-TextValue 'D14' '3.808.83'
$ws.Range('E14').Value = '  -0.59%  '
Set-TextValue 'D15' '8.39'
$ws.Range('E15').Value = '  -2.68%  '
Set-TextValue 'D16' '611.11'
$ws.Range('E16').Value = '  -3.26%  '
Set-TextValue 'D17' '65.585.30'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('E18').Value = '  +0.29%  '
Set-TextValue 'D19' '17.76'
$ws.Range('E19').Value = '  -1.87%  '
Set-TextValue 'D20' '3.283.53'
$ws.Range('E20').Value = '  -0.36%  '
Set-TextValue 'D21' '10.87'
$ws.Range('E21').Value = '  -4.11%  '
Set-TextValue 'D22' '0.886'
$ws.Range('E22').Value = '  -2.07%  '
Set-TextValue 'D23' '17.93'
$ws.Range('E23').Value = '  +0.52%  '
Set-TextValue 'D24' '100.42'
$ws.Range('E24').Value = '  -2.03%  '
Set-TextValue 'D25' '4.94'
$ws.Range('E25').Value = '  -0.38%  '
Set-TextValue 'D26' '4.01'
$ws.Range('E26').Value = '  +1.74%  '
Set-TextValue 'D27' '2.69'
$ws.Range('E27').Value = '  -0.61%  '
Set-TextValue 'D28' '9.41'
$ws.Range('E28').Value = '  -1.21%  '
Set-TextValue 'D29' '30.68'
$ws.Range('E29').Value = '  +1.75%  '
Set-TextValue 'D30' '8.39'
$ws.Range('E30').Value = '  -2.80%  '
Set-TextValue 'D31' '6.39'
$ws.Range('E31').Value = '  +0.56%  '
Set-TextValue 'D32' '3.69'
$ws.Range('E32').Value = '  -8.71%  '
Set-TextValue 'D33' '548.35'
$ws.Range('E33').Value = '  +0.77%  '
Set-TextValue 'D34' '10.81'
$ws.Range('E34').Value = '  -2.44%  '
Set-TextValue 'D35' '3.772.14'
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('E36').Value = '  -1.82%  '
Set-TextValue 'D37' '0.998'
$ws.Range('E37').Value = '  -0.15%  '
Set-TextValue 'D38' '55.92'
$ws.Range('E38').Value = '  -2.74%  '
$ws.Range('E39').Value = '  -0.85%  '
Set-TextValue 'D40' '32.27'
$ws.Range('E40').Value = '  -4.09%  '
$ws.Range('B41').Value = 'ApeXProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D41' '3.38'
$ws.Range('E41').Value = '  +3.47%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D42' '3.12'
$ws.Range('E42').Value = '  -3.72%  '
Set-TextValue 'D43' '2.57'
$ws.Range('E43').Value = '  -4.49%  '
Set-TextValue 'D44' '0.0₃0673'
$ws.Range('E44').Value = '  -8.68%  '
Set-TextValue 'D45' '0.329'
$ws.Range('E45').Value = '  -1.45%  '
Set-TextValue 'D46' '0.0404'
$ws.Range('E46').Value = '  -2.97%  '
Set-TextValue 'D47' '3.01'
$ws.Range('E47').Value = '  -6.78%  '
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('E49').Value = '  -2.20%  '
Set-TextValue 'D50' '2.49'
$ws.Range('E50').Value = '  -4.24%  '
Set-TextValue 'D51' '128.67'
$ws.Range('E51').Value = '  +4.82%  '

$excel.CutCopyMode = $false